$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells stay as text (avoid Excel auto-numeric conversion)
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.859.96'
$ws.Range("E2").Value = '  -1.99%  '

$ws.Range("D3").Value = '1.826.08'
$ws.Range("E3").Value = '  -2.38%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '239.47'
$ws.Range("E5").Value = '  -1.71%  '

$ws.Range("D6").Value = '0.6866'
$ws.Range("E6").Value = '  -2.77%  '

$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '0.07609'
$ws.Range("E8").Value = '  -3.34%  '

$ws.Range("D9").Value = '0.3015'
$ws.Range("E9").Value = '  -4.31%  '

$ws.Range("D10").Value = '23.43'
$ws.Range("E10").Value = '  -4.91%  '

$ws.Range("D11").Value = '0.07727'
$ws.Range("E11").Value = '  -3.50%  '

$ws.Range("D12").Value = '1.830.80'
$ws.Range("E12").Value = '  -2.72%  '

$ws.Range("D13").Value = '5.036'
$ws.Range("E13").Value = '  -3.37%  '

$ws.Range("D14").Value = '89.96'
$ws.Range("E14").Value = '  -4.39%  '

$ws.Range("E15").Value = '  -4.80%  '

$ws.Range("D16").Value = '6.397'
$ws.Range("E16").Value = '  -1.43%  '

$ws.Range("D17").Value = '0.000008255'
$ws.Range("E17").Value = '  -1.41%  '

$ws.Range("D18").Value = '28.861.88'
$ws.Range("E18").Value = '  -2.21%  '

$ws.Range("D19").Value = '242.28'
$ws.Range("E19").Value = '  -5.43%  '

$ws.Range("D20").Value = '2.088.45'
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("D21").Value = '12.59'
$ws.Range("E21").Value = '  -4.51%  '

$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").Value = '7.381'
$ws.Range("E23").Value = '  -3.12%  '

$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("E25").Value = '  -5.72%  '

$ws.Range("D26").Value = '160.09'
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("D27").Value = '8.682'
$ws.Range("E27").Value = '  -4.16%  '

$ws.Range("D28").Value = '18.09'
$ws.Range("E28").Value = '  -3.91%  '

$ws.Range("D29").Value = '1.529'
$ws.Range("E29").Value = '  +1.98%  '

$ws.Range("D30").Value = '4.190'
$ws.Range("E30").Value = '  -3.33%  '

$ws.Range("D31").Value = '4.140'
$ws.Range("E31").Value = '  -2.54%  '

$ws.Range("D32").Value = '1.190'
$ws.Range("E32").Value = '  -1.53%  '

$ws.Range("D33").Value = '0.05088'
$ws.Range("E33").Value = '  -4.42%  '

$ws.Range("D34").Value = '0.7555'
$ws.Range("E34").Value = '  +1.04%  '

$ws.Range("D35").Value = '1.808'
$ws.Range("E35").Value = '  -4.88%  '

$ws.Range("D36").Value = '1.137'
$ws.Range("E36").Value = '  -2.83%  '

$ws.Range("D37").Value = '2.689'
$ws.Range("E37").Value = '  -0.97%  '

$ws.Range("D38").Value = '0.01829'
$ws.Range("E38").Value = '  -2.47%  '

$ws.Range("D39").Value = '1.208.57'
$ws.Range("E39").Value = '  -4.51%  '

$ws.Range("D40").Value = '2.675'
$ws.Range("E40").Value = '  -2.71%  '

$ws.Range("D41").Value = '0.9106'
$ws.Range("E41").Value = '  +1.22%  '

$ws.Range("D42").Value = '108.44'
$ws.Range("E42").Value = '  -0.60%  '

$ws.Range("D43").Value = '0.9996'
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("D44").Value = '1.989.63'
$ws.Range("E44").Value = '  -2.38%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000122'
$ws.Range("E45").Value = '  -5.66%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.5160'
$ws.Range("E46").Value = '  -0.65%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '9.405'
$ws.Range("E47").Value = '  -1.09%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '5.308'
$ws.Range("E48").Value = '  -10.78%  '

$ws.Range("D49").Value = '62.31'
$ws.Range("E49").Value = '  -13.22%  '

$ws.Range("E50").Value = '  -5.38%  '

$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '0.4154'
$ws.Range("E51").Value = '  -3.90%  '
